$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 228; this shifts the existing rows 228-315
# down to 229-316 and extends the sheet dimension to A1:R316 automatically.
$ws.Rows("228:228").Insert()

# Populate the newly inserted row 228 with the new record data.
$ws.Cells.Item(228, 1).Value = 9
$ws.Cells.Item(228, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(228, 3).Value = "Metropolitana"
$ws.Cells.Item(228, 4).Value = 44795
$ws.Cells.Item(228, 5).Value = 13
$ws.Cells.Item(228, 6).Value = 100112043
$ws.Cells.Item(228, 7).Value = "Pepino ensalada"
$ws.Cells.Item(228, 8).Value = "Sin especificar"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 70
$ws.Cells.Item(228, 11).Value = 21000
$ws.Cells.Item(228, 12).Value = 22000
$ws.Cells.Item(228, 13).Value = 21500
$ws.Cells.Item(228, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(228, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(228, 16).Value = 358
$ws.Cells.Item(228, 17).Value = 60
$ws.Cells.Item(228, 18).Value = "Hortaliza"
